# Update the "give-me-some-credit" results table.
# The commit fixes handling of None metric values and computes metrics for
# results that were previously left un-optimized: row 2 (KAOGExp) gets new
# CERScore values, and the remaining method rows (3-14) are re-ordered to
# their corrected position while keeping their original metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final state for rows 2..14: Name, Success_Rate, CERScore_Distance_1..4
$rows = @(
    @{Row=2;  A="KAOGExp";      B=1;    C=9.06;               D=0.2420452041110067; E=0.06782578315039454; F=0.1025827849439245}
    @{Row=3;  A="cruds";        B=1;    C=10.48;              D=1.825658650695374;  E=0.8609219060674071;  F=0.673200966145333}
    @{Row=4;  A="wachter";      B=1;    C=10.22;              D=0.9682642677827201; E=0.5317357041877157;  F=0.5085406824339784}
    @{Row=5;  A="face-knn";     B=1;    C=9.19;               D=2.181617147650578;  E=1.148180030638227;   F=0.7671678082191781}
    @{Row=6;  A="revise";       B=1;    C=10.48;              D=1.830553903496672;  E=0.8574439824560099;  F=0.6735965235084805}
    @{Row=7;  A="cem";          B=1;    C=6.98;               D=0.9120619179605925; E=0.6591150092133984;  F=0.650797647921741}
    @{Row=8;  A="cem-vae";      B=1;    C=6.97;               D=0.9034640582713611; E=0.6568001494721392;  F=0.6500466030523181}
    @{Row=9;  A="dice";         B=1;    C=2.09;               D=1.329864808647734;  E=1.070182068133659;   F=0.827379661192294}
    @{Row=10; A="face-epsilon"; B=1;    C=9.039999999999999;  D=2.15136036891069;   E=1.161212482567081;   F=0.7563458904109589}
    @{Row=11; A="clue";         B=0.93; C=9.77;               D=1.638284066738282;  E=0.7815670584071454;  F=0.6477656169939259}
    @{Row=12; A="ar";           B=0.05; C=0.08;               D=0.02718818073701842;E=0.02194267304219456; F=0.02211416666666667}
    @{Row=13; A="cchvae";       B=0.93; C=9.77;               D=1.727530965816696;  E=0.8274086978646032;  F=0.6459519351013515}
    @{Row=14; A="gs";           B=1;    C=9.5;                D=0.6583560442159176; E=0.5082221524252166;  F=0.5198055902705477}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
